# Fix loading tags from excel
# Adds a "tags" column (with a "testing" sample value) to the
# physical_links and l3_links sheets, then leaves the workbook
# positioned the way the author left it when they saved: selection
# sitting in the new column, and "l3_links" as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- physical_links (sheet2) ---------------------------------------
$wsPhysicalLinks = $wb.Worksheets.Item("physical_links")
$wsPhysicalLinks.Range("L1").Value = "tags"
$wsPhysicalLinks.Range("L2").Value = "testing"
[void]$wsPhysicalLinks.Range("L2").Select()

# --- l3_links (sheet3) ----------------------------------------------
$wsL3Links = $wb.Worksheets.Item("l3_links")
$wsL3Links.Range("P1").Value = "tags"
$wsL3Links.Range("P2").Value = "testing"

# l3_links ends up as the active sheet/tab, selection resting one row
# below the newly entered data (P3).
[void]$wsL3Links.Activate()
[void]$wsL3Links.Range("P3").Select()
